# Applies the recorded change set:
#  - Insert a new column D ("weightKeep2") into the "generate-weights" sheet,
#    pushing the old D column (Pos-rate) to column E.
#  - Fill the new D column with header + 5 computed values, using a
#    wrap-text style.
#  - Update the active sheet / selection bookkeeping: "generate-weights"
#    (4th tab) becomes the active sheet/tab, "generate-weights_test"
#    (1st tab) is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item(1)
$wsGen  = $wb.Worksheets.Item(4)

# --- "generate-weights" sheet: insert new column D ---------------------
$wsGen.Columns.Item(4).Insert()

$wsGen.Range("D1").Value = "weightKeep2"
$wsGen.Range("D2").Value = 0.0975373359203966
$wsGen.Range("D3").Value = 0.426216532321556
$wsGen.Range("D4").Value = 0.636270697849109
$wsGen.Range("D5").Value = 0.818623193631408
$wsGen.Range("D6").Value = 0.839010519465712

$wsGen.Range("D2:D6").WrapText = $true

# --- selection / active-sheet bookkeeping -------------------------------
$null = $wsTest.Range("C2").Select()

$null = $wsGen.Activate()
$null = $wsGen.Range("D2").Select()
